# Merged the tests and updated the excel
#
# - Disable concurrent/multi-threaded calculation (workbook-level calc setting,
#   xlsx: calcPr/@concurrentCalc="0").
# - Change the two "Runmode" cells that were still "Y" (C6:C7) to "N", adding
#   a new shared string for it.
# - Move the active selection from C4 to C7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Workbook calculation option mirroring calcPr concurrentCalc="0"
$excel.MultiThreadedCalculation.Enabled = $false

# Data edits: last two rows' Runmode flips from Y to N
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"

# Selection moves to C7
$ws.Range("C7").Select()
